$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily COVID data rows (day, month, year, Pruebas Realizadas, Pruebas Positivas)
$data = @(
    @(11, 1, 2021, 5899, 1567),
    @(12, 1, 2021, 5751, 1004),
    @(13, 1, 2021, 6479, 951),
    @(14, 1, 2021, 1175, 623),
    @(15, 1, 2021, 5670, 1038),
    @(16, 1, 2021, 2221, 290)
)

$startRow = 290
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $prev = $r - 1

    # Copy format (style + row height) from the previous row's A:D cells
    $ws.Range("A" + $prev + ":D" + $prev).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 15

    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 1).Formula = "=+Hoja1!`$B" + $r + "&""/""&Hoja1!`$C" + $r + "&""/""&Hoja1!`$D" + $r
}

# Grow the table to include the newly added rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:I295"))

# Match the saved selection/view reported in the workbook after the edit
$ws.Range("E296").Select()

Write-Output "OK"
